$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 690.8333
$ws.Range("I32").Value = 690.8333
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 690.8333
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -364.8333
$ws.Range("N32").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 335198.66
$ws.Range("I125").Value = 600
$ws.Range("J125").Value = 502498
$ws.Range("K125").Value = 5400
$ws.Range("L125").Value = 4522482
$ws.Range("M125").Value = -2940
$ws.Range("N125").Value = -4527402

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 780000
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 780000
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 780000
$ws.Range("N134").Value = -790140

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2814.8333
$ws.Range("I137").Value = 2814.8333
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 8444.499899999999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -5894.499899999999
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5297.5
$ws.Range("I32").Value = 1353.5652
$ws.Range("J32").Value = 18256.143
$ws.Range("K32").Value = 1353.5652
$ws.Range("L32").Value = 18256.143
$ws.Range("M32").Value = -1066.5652
$ws.Range("N32").Value = -18830.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 20000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 20000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20976

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1567.5454
$ws.Range("I61").Value = 1567.5454
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1567.5454
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1355.5454
$ws.Range("N61").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1314
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6568
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1798.3334
$ws.Range("I74").Value = 1798.3334
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1798.3334
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -924.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1798.3334
$ws.Range("I77").Value = 1798.3334
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 8991.666999999999
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -4623.666999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 29999
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 29999
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 29999
$ws.Range("N80").Value = -31995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 29999
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 29999
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 89997
$ws.Range("N83").Value = -99981

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1567.5454
$ws.Range("I136").Value = 1567.5454
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4702.6362
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2152.6362
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3459.6
$ws.Range("I20").Value = 3132.6667
$ws.Range("J20").Value = 3950
$ws.Range("K20").Value = 3132.6667
$ws.Range("L20").Value = 3950
$ws.Range("M20").Value = -2885.6667
$ws.Range("N20").Value = -4444

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 14996.875
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 14996.875
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 14996.875
$ws.Range("N35").Value = -15616.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 21573.588
$ws.Range("I82").Value = 15675.6
$ws.Range("J82").Value = 29999.285
$ws.Range("K82").Value = 15675.6
$ws.Range("L82").Value = 29999.285
$ws.Range("M82").Value = -15292.6
$ws.Range("N82").Value = -30765.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 21573.588
$ws.Range("I85").Value = 15675.6
$ws.Range("J85").Value = 29999.285
$ws.Range("K85").Value = 15675.6
$ws.Range("L85").Value = 29999.285
$ws.Range("M85").Value = -14349.6
$ws.Range("N85").Value = -32651.285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1764.2
$ws.Range("I94").Value = 1193.4667
$ws.Range("J94").Value = 3476.4
$ws.Range("K94").Value = 1193.4667
$ws.Range("L94").Value = 3476.4
$ws.Range("M94").Value = -742.4666999999999
$ws.Range("N94").Value = -4378.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 983
$ws.Range("I107").Value = 950
$ws.Range("J107").Value = 999.5
$ws.Range("K107").Value = 950
$ws.Range("L107").Value = 999.5
$ws.Range("M107").Value = 970
$ws.Range("N107").Value = -4839.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7224.722
$ws.Range("I134").Value = 6561.4707
$ws.Range("J134").Value = 18500
$ws.Range("K134").Value = 19684.4121
$ws.Range("L134").Value = 55500
$ws.Range("M134").Value = -17149.4121
$ws.Range("N134").Value = -60570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 29448.4
$ws.Range("I59").Value = 24504
$ws.Range("J59").Value = 29997.777
$ws.Range("K59").Value = 24504
$ws.Range("L59").Value = 29997.777
$ws.Range("M59").Value = -23359
$ws.Range("N59").Value = -32287.777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5499.75
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 5999.6665
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 5999.6665
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -7247.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5499.75
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 5999.6665
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 29998.3325
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -36238.3325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 648.4
$ws.Range("I105").Value = 640.7143
$ws.Range("J105").Value = 666.3333
$ws.Range("K105").Value = 640.7143
$ws.Range("L105").Value = 666.3333
$ws.Range("M105").Value = 1106.2857
$ws.Range("N105").Value = -4160.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 619.82355
$ws.Range("I5").Value = 569.13336
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1707.40008
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1595.40008
$ws.Range("N5").Value = -3224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 619.82355
$ws.Range("I135").Value = 569.13336
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 5122.20024
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -2587.20024
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 604000
$ws.Range("I7").Value = 666666.7
$ws.Range("J7").Value = 510000
$ws.Range("K7").Value = 666666.7
$ws.Range("L7").Value = 510000
$ws.Range("M7").Value = -666554.7
$ws.Range("N7").Value = -510224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 604000
$ws.Range("I8").Value = 666666.7
$ws.Range("J8").Value = 510000
$ws.Range("K8").Value = 666666.7
$ws.Range("L8").Value = 510000
$ws.Range("M8").Value = -666527.7
$ws.Range("N8").Value = -510278

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 8000
$ws.Range("I46").Value = 6666.6665
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 6666.6665
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = -6510.6665
$ws.Range("N46").Value = -10312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 18750
$ws.Range("I57").Value = 15000
$ws.Range("J57").Value = 20000
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 20000
$ws.Range("M57").Value = -14180
$ws.Range("N57").Value = -21640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2844.2727
$ws.Range("I122").Value = 2224.875
$ws.Range("J122").Value = 3198.2144
$ws.Range("K122").Value = 6674.625
$ws.Range("L122").Value = 9594.643199999999
$ws.Range("M122").Value = -4224.625
$ws.Range("N122").Value = -14494.6432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 25749
$ws.Range("I68").Value = 7499.6665
$ws.Range("J68").Value = 36698.6
$ws.Range("K68").Value = 7499.6665
$ws.Range("L68").Value = 36698.6
$ws.Range("M68").Value = -6750.6665
$ws.Range("N68").Value = -38196.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 25749
$ws.Range("I71").Value = 7499.6665
$ws.Range("J71").Value = 36698.6
$ws.Range("K71").Value = 37498.3325
$ws.Range("L71").Value = 183493
$ws.Range("M71").Value = -33754.3325
$ws.Range("N71").Value = -190981

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3323.1052
$ws.Range("I136").Value = 2928.5
$ws.Range("J136").Value = 3999.5715
$ws.Range("K136").Value = 8785.5
$ws.Range("L136").Value = 11998.7145
$ws.Range("M136").Value = -6235.5
$ws.Range("N136").Value = -17098.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 6225
$ws.Range("I14").Value = 2400
$ws.Range("J14").Value = 7500
$ws.Range("K14").Value = 2400
$ws.Range("L14").Value = 7500
$ws.Range("M14").Value = -2232
$ws.Range("N14").Value = -7836

Write-Output "Applied all changes"